$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra trailing row (8: formate / engineered) and the
# CO2 / engineered column (F) that only existed to support it.
$ws.Rows.Item(8).Delete()
$ws.Columns.Item(6).Delete()

# Reorder the data rows: "nucleic acids" (currently row 5) moves up to
# row 3, pushing "amino acids" (row 3) down to row 4 and "sugars"
# (row 4) down to row 5. Rewrite rows 3-5 directly with their final
# contents instead of relying on cut/insert semantics.
$ws.Cells.Item(3, 1).Value = "nucleic acids"
$ws.Cells.Item(3, 2).Value = "experiment"
$ws.Cells.Item(3, 3).Value = "experiment"
$ws.Cells.Item(3, 4).Value = "model"
$ws.Cells.Item(3, 5).ClearContents()

$ws.Cells.Item(4, 1).Value = "amino acids"
$ws.Cells.Item(4, 2).Value = "experiment"
$ws.Cells.Item(4, 3).Value = "experiment"
$ws.Cells.Item(4, 4).Value = "model"
$ws.Cells.Item(4, 5).ClearContents()

$ws.Cells.Item(5, 1).Value = "sugars"
$ws.Cells.Item(5, 2).Value = "experiment"
$ws.Cells.Item(5, 3).Value = "experiment"
$ws.Cells.Item(5, 4).Value = "experiment"
$ws.Cells.Item(5, 5).Value = "experiment"

# Rows 6 (fatty acids) and 7 (aromatics) are unchanged in content.

# Match the saved selection from the authored workbook.
[void]$ws.Range("E6").Select()
